$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for this product/market. It belongs
# chronologically at row 9 (Fecha 44469), so insert a fresh row there and
# push the existing rows 9-16 down to 10-17 (matches the diff: dimension
# grows from A1:R16 to A1:R17).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new record's data. Columns
# A, B, C, E, F, G, H, N, Q, R are constant for every row in this sheet.
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 44469
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 100112045
$ws.Range("G9").Value = "Zapallo"
$ws.Range("H9").Value = "Camote"
$ws.Range("I9").Value = "1a nueva(o)"
$ws.Range("J9").Value = 1200
$ws.Range("K9").Value = 600
$ws.Range("L9").Value = 650
$ws.Range("M9").Value = 625
$ws.Range("N9").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O9").Value = "Perú"
$ws.Range("P9").Value = 625
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = "Hortaliza"
